# Updated cryptos list on Wed May 15 07:24:19 UTC 2024 with GitHub Actions
#
# The "Price" column (D) holds numeric-looking strings (e.g. "568.18",
# "0.500") that must stay as TEXT, matching the source feed's inlineStr
# cells. Plain `.Value = "..."` lets Excel auto-coerce those into real
# numbers (losing trailing zeros / exact formatting), so each D-column
# write is wrapped: force text format, assign, then restore the default
# style so no stray style index is left on the cell.
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "61.932.94"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.906.04"
$ws.Range("E3").Value = "  -0.39%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "568.18"
$ws.Range("E5").Value = "  -3.38%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "143.57"
$ws.Range("E6").Value = "  -1.88%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.11%  "

# Rows 8 & 9 swap: LidoStakedEther and XRP trade ranking positions
Set-TextValue $ws "B8" "LidoStakedEther"
Set-TextValue $ws "C8" "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue $ws "D8" "2.903.16"
$ws.Range("E8").Value = "  -0.49%  "

Set-TextValue $ws "B9" "XRP"
Set-TextValue $ws "C9" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws "D9" "0.500"
$ws.Range("E9").Value = "  -0.99%  "

# Row 10 - Toncoin
Set-TextValue $ws "D10" "6.96"
$ws.Range("E10").Value = "  -0.76%  "

# Row 11 - Dogecoin
Set-TextValue $ws "D11" "0.148"
$ws.Range("E11").Value = "  -1.98%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.49%  "

# Row 13 - ShibaInu
Set-TextValue $ws "D13" "0.0000236"
$ws.Range("E13").Value = "  -0.15%  "

# Row 14 - Avalanche
Set-TextValue $ws "D14" "32.48"
$ws.Range("E14").Value = "  +0.14%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.17%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D16" "3.387.90"
$ws.Range("E16").Value = "  -0.42%  "

# Row 17 - WrappedBTC
Set-TextValue $ws "D17" "61.883.75"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18 - WrappedEther
Set-TextValue $ws "D18" "2.900.69"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.66%  "

# Row 20 - BitcoinCash
Set-TextValue $ws "D20" "429.32"
$ws.Range("E20").Value = "  -1.31%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  -2.75%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -1.40%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -1.08%  "

# Row 24 - Litecoin
Set-TextValue $ws "D24" "78.79"
$ws.Range("E24").Value = "  -2.38%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue $ws "D25" "11.97"
$ws.Range("E25").Value = "  +0.74%  "

# Row 26 - RenderToken
$ws.Range("E26").Value = "  -7.78%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.10%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  -3.83%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +10.72%  "

# Row 30 - NEARProtocol
Set-TextValue $ws "D30" "6.95"
$ws.Range("E30").Value = "  -5.20%  "

# Row 31 - PancakeSwap
Set-TextValue $ws "D31" "2.51"
$ws.Range("E31").Value = "  -2.77%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -5.35%  "

# Row 33 - FirstDigitalUSD
Set-TextValue $ws "D33" "0.999"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -3.31%  "

# Row 35 - EthereumClassic
Set-TextValue $ws "D35" "25.58"
$ws.Range("E35").Value = "  -1.91%  "

# Row 36 - Mantle
Set-TextValue $ws "D36" "0.952"
$ws.Range("E36").Value = "  -2.62%  "

# Row 37 - Filecoin
Set-TextValue $ws "D37" "5.39"
$ws.Range("E37").Value = "  -2.68%  "

# Row 38 - OKB
Set-TextValue $ws "D38" "48.81"
$ws.Range("E38").Value = "  -0.87%  "

# Row 39 - dogwifhat
$ws.Range("E39").Value = "  -6.60%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -4.94%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -1.12%  "

# Row 42 - Arweave
Set-TextValue $ws "D42" "41.09"
$ws.Range("E42").Value = "  +5.16%  "

# Row 43 - Cosmos
Set-TextValue $ws "D43" "8.16"
$ws.Range("E43").Value = "  -2.54%  "

# Row 44 - TheGraph
$ws.Range("E44").Value = "  -3.13%  "

# Row 45 - Maker
Set-TextValue $ws "D45" "2.710.42"
$ws.Range("E45").Value = "  +0.47%  "

# Row 46 - Monero
Set-TextValue $ws "D46" "133.09"
$ws.Range("E46").Value = "  -1.45%  "

# Row 47 - VeChain
Set-TextValue $ws "D47" "0.0336"
$ws.Range("E47").Value = "  -0.59%  "

# Row 48 - Bittensor
Set-TextValue $ws "D48" "348.65"
$ws.Range("E48").Value = "  +0.14%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  -0.03%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  -1.05%  "

# Row 51 - FLOKI
Set-TextValue $ws "D51" "0.000209"
$ws.Range("E51").Value = "  +10.26%  "
